$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update indicator_last_checked timestamps for existing rows 2 and 3
$ws.Range("C2").Value = "2025-09-24T15:10:35.335491"
$ws.Range("C3").Value = "2025-09-24T15:10:37.042674"

# Fill in full Technical Indicators data for rows 4-11 (previously empty placeholder cells)
# Row 4: ABBV
$row4 = New-Object 'object[,]' 1,21
$row4[0,0] = "https://www.investing.com/equities/abbvie-inc-technical"
$row4[0,1] = "2025-09-24T15:10:38.971697"
$row4[0,2] = "mock"
$row4[0,3] = "Mock data used due to network unavailability"
$row4[0,4] = 354
$row4[0,5] = 345
$row4[0,6] = 325
$row4[0,7] = 379
$row4[0,8] = 399
$row4[0,9] = 359
$row4[0,10] = 362
$row4[0,11] = 32
$row4[0,12] = -8
$row4[0,13] = 4.5
$row4[0,14] = -2
$row4[0,15] = 384
$row4[0,16] = 362
$row4[0,17] = 340
$row4[0,18] = 2709762
$row4[0,19] = 32
$row4[0,20] = 3
$ws.Range("B4:V4").Value = $row4

# Row 5: ADBE
$row5 = New-Object 'object[,]' 1,21
$row5[0,0] = "https://www.investing.com/equities/adobe-inc-technical"
$row5[0,1] = "2025-09-24T15:10:40.686815"
$row5[0,2] = "mock"
$row5[0,3] = "Mock data used due to network unavailability"
$row5[0,4] = 176
$row5[0,5] = 165
$row5[0,6] = 155
$row5[0,7] = 201
$row5[0,8] = 211
$row5[0,9] = 181
$row5[0,10] = 180
$row5[0,11] = 53
$row5[0,12] = -7
$row5[0,13] = -4.5
$row5[0,14] = 3
$row5[0,15] = 206
$row5[0,16] = 183
$row5[0,17] = 160
$row5[0,18] = 8699583
$row5[0,19] = 53
$row5[0,20] = 4
$ws.Range("B5:V5").Value = $row5

# Row 6: ADI
$row6 = New-Object 'object[,]' 1,21
$row6[0,0] = "https://www.investing.com/equities/analog-devices-technical"
$row6[0,1] = "2025-09-24T15:10:42.673470"
$row6[0,2] = "mock"
$row6[0,3] = "Mock data used due to network unavailability"
$row6[0,4] = 274
$row6[0,5] = 255
$row6[0,6] = 240
$row6[0,7] = 289
$row6[0,8] = 304
$row6[0,9] = 269
$row6[0,10] = 270
$row6[0,11] = 62
$row6[0,12] = 2
$row6[0,13] = -0.5
$row6[0,14] = -4
$row6[0,15] = 294
$row6[0,16] = 272
$row6[0,17] = 250
$row6[0,18] = 10933672
$row6[0,19] = 42
$row6[0,20] = 3
$ws.Range("B6:V6").Value = $row6

# Row 7: ADSK
$row7 = New-Object 'object[,]' 1,21
$row7[0,0] = "https://www.investing.com/equities/autodesk-inc-technical"
$row7[0,1] = "2025-09-24T15:10:43.838381"
$row7[0,2] = "mock"
$row7[0,3] = "Mock data used due to network unavailability"
$row7[0,4] = 112
$row7[0,5] = 95
$row7[0,6] = 80
$row7[0,7] = 127
$row7[0,8] = 142
$row7[0,9] = 107
$row7[0,10] = 108
$row7[0,11] = 41
$row7[0,12] = 1
$row7[0,13] = -1.5
$row7[0,14] = -1
$row7[0,15] = 132
$row7[0,16] = 111
$row7[0,17] = 90
$row7[0,18] = 6281011
$row7[0,19] = 31
$row7[0,20] = 2
$ws.Range("B7:V7").Value = $row7

# Row 8: ADYEY
$row8 = New-Object 'object[,]' 1,21
$row8[0,0] = "https://www.investing.com/equities/adyen-nv-otc-technical"
$row8[0,1] = "2025-09-24T15:10:45.441807"
$row8[0,2] = "mock"
$row8[0,3] = "Mock data used due to network unavailability"
$row8[0,4] = 226
$row8[0,5] = 195
$row8[0,6] = 190
$row8[0,7] = 241
$row8[0,8] = 246
$row8[0,9] = 221
$row8[0,10] = 218
$row8[0,11] = 48
$row8[0,12] = 8
$row8[0,13] = -4.5
$row8[0,14] = -2
$row8[0,15] = 246
$row8[0,16] = 218
$row8[0,17] = 190
$row8[0,18] = 5834618
$row8[0,19] = 38
$row8[0,20] = 9
$ws.Range("B8:V8").Value = $row8

# Row 9: AEM
$row9 = New-Object 'object[,]' 1,21
$row9[0,0] = "https://www.investing.com/equities/agnico-eagle-mines-technical"
$row9[0,1] = "2025-09-24T15:10:47.092396"
$row9[0,2] = "mock"
$row9[0,3] = "Mock data used due to network unavailability"
$row9[0,4] = 402
$row9[0,5] = 385
$row9[0,6] = 375
$row9[0,7] = 427
$row9[0,8] = 437
$row9[0,9] = 407
$row9[0,10] = 406
$row9[0,11] = 36
$row9[0,12] = -4
$row9[0,13] = -1.5
$row9[0,14] = 2
$row9[0,15] = 432
$row9[0,16] = 406
$row9[0,17] = 380
$row9[0,18] = 8205806
$row9[0,19] = 26
$row9[0,20] = 7
$ws.Range("B9:V9").Value = $row9

# Row 10: AJG
$row10 = New-Object 'object[,]' 1,21
$row10[0,0] = "https://www.investing.com/equities/arthur-j.-gallagher---co-technical"
$row10[0,1] = "2025-09-24T15:10:48.502007"
$row10[0,2] = "mock"
$row10[0,3] = "Mock data used due to network unavailability"
$row10[0,4] = 450
$row10[0,5] = 435
$row10[0,6] = 420
$row10[0,7] = 465
$row10[0,8] = 480
$row10[0,9] = 445
$row10[0,10] = 446
$row10[0,11] = 60
$row10[0,12] = 0
$row10[0,13] = -2.5
$row10[0,14] = 2
$row10[0,15] = 470
$row10[0,16] = 450
$row10[0,17] = 430
$row10[0,18] = 4003350
$row10[0,19] = 20
$row10[0,20] = 1
$ws.Range("B10:V10").Value = $row10

# Row 11: ALL
$row11 = New-Object 'object[,]' 1,21
$row11[0,0] = "https://www.investing.com/equities/allstate-corp-technical"
$row11[0,1] = "2025-09-24T15:10:50.048240"
$row11[0,2] = "mock"
$row11[0,3] = "Mock data used due to network unavailability"
$row11[0,4] = 312
$row11[0,5] = 305
$row11[0,6] = 290
$row11[0,7] = 337
$row11[0,8] = 352
$row11[0,9] = 317
$row11[0,10] = 324
$row11[0,11] = 51
$row11[0,12] = -9
$row11[0,13] = -1.5
$row11[0,14] = 1
$row11[0,15] = 342
$row11[0,16] = 321
$row11[0,17] = 300
$row11[0,18] = 6499221
$row11[0,19] = 41
$row11[0,20] = 2
$ws.Range("B11:V11").Value = $row11
